{"js": "// Replace the three-digit x one-digit multiplication problems throughout\n// the document's table with a new set of problems/answers, per the\n// commit's updated output. Each old expression is unique in the\n// document, so a simple text search-and-replace per pair is safe.\nconst replacements = [\n  [\"164\u00d79=1476\", \"322\u00d74=1288\"],\n  [\"790\u00d75=3950\", \"997\u00d75=4985\"],\n  [\"905\u00d74=3620\", \"590\u00d79=5310\"],\n  [\"271\u00d77=1897\", \"544\u00d75=2720\"],\n  [\"596\u00d78=4768\", \"529\u00d79=4761\"],\n  [\"488\u00d73=1464\", \"952\u00d78=7616\"],\n  [\"354\u00d78=2832\", \"403\u00d77=2821\"],\n  [\"234\u00d72=468\", \"955\u00d74=3820\"],\n  [\"933\u00d75=4665\", \"102\u00d77=714\"],\n  [\"943\u00d73=2829\", \"909\u00d73=2727\"],\n  [\"317\u00d77=2219\", \"606\u00d76=3636\"],\n  [\"978\u00d73=2934\", \"168\u00d78=1344\"],\n  [\"511\u00d78=4088\", \"654\u00d73=1962\"],\n  [\"373\u00d76=2238\", \"613\u00d77=4291\"],\n  [\"290\u00d72=580\", \"826\u00d75=4130\"],\n  [\"452\u00d75=2260\", \"744\u00d79=6696\"],\n  [\"304\u00d73=912\", \"452\u00d73=1356\"],\n  [\"534\u00d76=3204\", \"212\u00d76=1272\"],\n  [\"794\u00d78=6352\", \"108\u00d77=756\"],\n  [\"115\u00d74=460\", \"162\u00d79=1458\"],\n  [\"924\u00d73=2772\", \"723\u00d74=2892\"],\n  [\"362\u00d73=1086\", \"695\u00d77=4865\"],\n  [\"366\u00d76=2196\", \"806\u00d76=4836\"],\n  [\"998\u00d78=7984\", \"659\u00d72=1318\"],\n  [\"234\u00d75=1170\", \"530\u00d72=1060\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication problems throughout\n# the document's table with a new set of problems/answers, per the\n# commit's updated output. Each old expression is unique in the\n# document, so Find/Execute with exact text match is safe for each pair.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"164\u00d79=1476\", \"322\u00d74=1288\"),\n    @(\"790\u00d75=3950\", \"997\u00d75=4985\"),\n    @(\"905\u00d74=3620\", \"590\u00d79=5310\"),\n    @(\"271\u00d77=1897\", \"544\u00d75=2720\"),\n    @(\"596\u00d78=4768\", \"529\u00d79=4761\"),\n    @(\"488\u00d73=1464\", \"952\u00d78=7616\"),\n    @(\"354\u00d78=2832\", \"403\u00d77=2821\"),\n    @(\"234\u00d72=468\", \"955\u00d74=3820\"),\n    @(\"933\u00d75=4665\", \"102\u00d77=714\"),\n    @(\"943\u00d73=2829\", \"909\u00d73=2727\"),\n    @(\"317\u00d77=2219\", \"606\u00d76=3636\"),\n    @(\"978\u00d73=2934\", \"168\u00d78=1344\"),\n    @(\"511\u00d78=4088\", \"654\u00d73=1962\"),\n    @(\"373\u00d76=2238\", \"613\u00d77=4291\"),\n    @(\"290\u00d72=580\", \"826\u00d75=4130\"),\n    @(\"452\u00d75=2260\", \"744\u00d79=6696\"),\n    @(\"304\u00d73=912\", \"452\u00d73=1356\"),\n    @(\"534\u00d76=3204\", \"212\u00d76=1272\"),\n    @(\"794\u00d78=6352\", \"108\u00d77=756\"),\n    @(\"115\u00d74=460\", \"162\u00d79=1458\"),\n    @(\"924\u00d73=2772\", \"723\u00d74=2892\"),\n    @(\"362\u00d73=1086\", \"695\u00d77=4865\"),\n    @(\"366\u00d76=2196\", \"806\u00d76=4836\"),\n    @(\"998\u00d78=7984\", \"659\u00d72=1318\"),\n    @(\"234\u00d75=1170\", \"530\u00d72=1060\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
